$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the value to be stored as text even when it looks numeric
    # (e.g. "208.40"), mirroring Excel's "number stored as text" behavior,
    # then reset the style so no extra formatting sticks to the cell.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.701.44"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.592.29"
$ws.Range("E3").Value = "  -2.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.17%  "

# Row 5 - BNB
Set-TextValue "D5" "208.40"
$ws.Range("E5").Value = "  -1.65%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.88%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.24%  "

# Row 8 - Solana
Set-TextValue "D8" "22.29"
$ws.Range("E8").Value = "  -4.11%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.98%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.47%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0867"
$ws.Range("E11").Value = "  -1.56%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.818.27"
$ws.Range("E12").Value = "  -2.35%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.595.50"
$ws.Range("E13").Value = "  -2.51%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -3.78%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -4.32%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "27.684.34"
$ws.Range("E16").Value = "  -0.87%  "

# Row 17 - Litecoin
Set-TextValue "D17" "63.42"
$ws.Range("E17").Value = "  -2.20%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "220.19"
$ws.Range("E18").Value = "  -3.44%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  -3.01%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.34"
$ws.Range("E20").Value = "  -3.79%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -4.54%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.67"
$ws.Range("E23").Value = "  -3.10%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -3.89%  "

# Row 25 - Monero
Set-TextValue "D25" "154.03"
$ws.Range("E25").Value = "  -0.40%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -1.90%  "

# Row 27 - BinanceUSD
$ws.Range("E27").Value = "  +0.23%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.15"
$ws.Range("E28").Value = "  -1.53%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -4.89%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.46%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -5.16%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.374.00"
$ws.Range("E33").Value = "  -3.22%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -4.89%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -5.00%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -2.76%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.25%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -1.20%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -3.19%  "

# Row 40 - ARBITRUM
Set-TextValue "D40" "0.828"
$ws.Range("E40").Value = "  -2.77%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.23%  "

# Row 42 - WEMIXToken
Set-TextValue "D42" "0.967"
$ws.Range("E42").Value = "  -3.74%  "

# Row 43 - Aave
Set-TextValue "D43" "64.64"
$ws.Range("E43").Value = "  -1.74%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +2.42%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -3.66%  "

# Rows 46-47: RenderToken and RocketPoolETH swap order
# Row 46 becomes RocketPoolETH, Row 47 becomes RenderToken
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.729.12"
$ws.Range("E46").Value = "  -2.36%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "1.73"
$ws.Range("E47").Value = "  -5.26%  "

# Row 48 - Quant
Set-TextValue "D48" "87.21"
$ws.Range("E48").Value = "  -1.64%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -0.77%  "

# Row 50 - Algorand
Set-TextValue "D50" "0.0966"
$ws.Range("E50").Value = "  -4.04%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -1.63%  "
